# SWAN.xlsx update
# The workbook lists the 18 SWAN questions (rows 3-20) plus 3 scoring rows
# (22-24). Column E ("Value Labels") already carries the full -3..3 value
# scale text (and its formatting) on row 3 (E3); every other data row is
# missing that text. This script copies the value + formatting from E3
# down onto every other question/scoring row, narrows column E so the
# text wraps, and leaves the selection on the last edited block (E22:E24)
# the way the author's Excel session ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full text of the value-scale note already stored in E3 (shared string).
$valueScaleText = "-3= Far above average`n-2= Above average`n-1= Slightly above average`n0= Average`n1= Slightly below average`n2= Above average`n3= Far above average"

# Rows that need the value-scale text in column E (E3 already has it).
$targetRows = @(4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,22,23,24)

# 1) Write the text first so the subsequent format-only paste (below)
#    doesn't need to invent a brand new style for "value + border + wrap".
foreach ($r in $targetRows) {
    $ws.Cells.Item($r, 5).Value2 = $valueScaleText
}

# 2) Copy E3's formatting (quote-prefix number format, bottom border,
#    wrapped text) onto each of those cells so they all match E3 exactly.
$ws.Range("E3").Copy()
foreach ($r in $targetRows) {
    $ws.Cells.Item($r, 5).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

# 3) Column E becomes a fixed, narrower width now that it holds wrapped
#    multi-line text instead of a single long best-fit line.
$ws.Columns.Item(5).ColumnWidth = 37

# 4) Rows now auto-size taller to show the wrapped 7-line value scale.
$dataRowHeight = 109.2
foreach ($r in @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,22,23,24)) {
    $ws.Rows.Item($r).RowHeight = $dataRowHeight
}
# Row 20 has the thick bottom border and ends up very slightly taller.
$ws.Rows.Item(20).RowHeight = 109.8

# 5) Leave the selection where the author's session ended up.
$ws.Range("E22:E24").Select()
